# ISPAuthorization Provider and creation output table
# The "Variable" sheet lists indicator variables IND_1..IND_14 used to build
# a SNDG-based output table. This renumbers the first nine indicators
# (IND_1..IND_9 -> IND_01..IND_09, zero-padded to match IND_10..IND_14) and
# moves the selection to reflect the next empty row for further entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variable")

$ws.Range("E4").Value  = "IND_01                         "
$ws.Range("E5").Value  = "IND_02"
$ws.Range("E6").Value  = "IND_03"
$ws.Range("E7").Value  = "IND_04"
$ws.Range("E8").Value  = "IND_05"
$ws.Range("E9").Value  = "IND_06"
$ws.Range("E10").Value = "IND_07"
$ws.Range("E11").Value = "IND_08"
$ws.Range("E12").Value = "IND_09"

$ws.Range("E20").Select()
